$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for columns G:T, rows 2-26
$data = @{
    2 = @(421.7027586666666, 1265.108276, 0.8230205889389596, 0.8230205889389596, 3, 1, 43.19793166666667, 129.593795, 0.7412538312889448, 0.7412538312889448, 18216.68695252749, 163950.1825727474, 0.6100671647806876, 0.6100671647806876)
    3 = @(421.7027586666666, 1265.108276, 0.8230205889389596, 0.8230205889389596, 3, 1, 10.38032666666667, 31.14098, 0.1781209566020688, 0.1781209566020688, 4377.412391194497, 39396.71152075048, 0.1465972146050055, 0.1465972146050055)
    4 = @(421.7027586666666, 1265.108276, 0.8230205889389596, 0.8230205889389596, 2, 0.6666666666666666, 0.422089, 1.266267, 0.007242825670663926, 0.007242825670663927, 177.9960957028546, 1601.964861325692, 0.005960994649052039, 0.00596099464905204)
    5 = @(421.7027586666666, 1265.108276, 0.8230205889389596, 0.8230205889389596, 3, 1, 4.133026333333333, 12.399079, 0.07092056230936286, 0.07092056230936288, 1742.908606408644, 15686.1774576778, 0.058369082959734, 0.05836908295973401)
    6 = @(421.7027586666666, 1265.108276, 0.8230205889389596, 0.8230205889389596, 3, 1, 0.1434673333333333, 0.430402, 0.002461824128959449, 0.002461824128959449, 60.50057024521688, 544.505132206952, 0.002026131944480347, 0.002026131944480347)
    7 = @(69.70494733333334, 209.114842, 0.136040387754698, 0.136040387754698, 3, 1, 43.19793166666667, 129.593795, 0.7412538312889448, 0.7412538312889448, 3011.109551733933, 27099.98596560539, 0.1008404586332036, 0.1008404586332036)
    8 = @(69.70494733333334, 209.114842, 0.136040387754698, 0.136040387754698, 3, 1, 10.38032666666667, 31.14098, 0.1781209566020688, 0.1781209566020688, 723.5601236027956, 6512.04111242516, 0.02423164400338318, 0.02423164400338318)
    9 = @(69.70494733333334, 209.114842, 0.136040387754698, 0.136040387754698, 2, 0.6666666666666666, 0.422089, 1.266267, 0.007242825670663926, 0.007242825670663927, 29.42169151497933, 264.795223634814, 0.0009853168126768013, 0.0009853168126768013)
    10 = @(69.70494733333334, 209.114842, 0.136040387754698, 0.136040387754698, 3, 1, 4.133026333333333, 12.399079, 0.07092056230936286, 0.07092056230936288, 288.0923828922797, 2592.831446030518, 0.009648060796346945, 0.009648060796346947)
    11 = @(69.70494733333334, 209.114842, 0.136040387754698, 0.136040387754698, 3, 1, 0.1434673333333333, 0.430402, 0.002461824128959449, 0.002461824128959449, 10.00038291405378, 90.003446226484, 0.0003349075090875152, 0.0003349075090875152)
    12 = @(1.088159666666667, 3.264479, 0.002123718167154624, 0.002123718167154624, 3, 1, 43.19793166666667, 129.593795, 0.7412538312889448, 0.7412538312889448, 47.00624692308944, 423.056222307805, 0.001574214227981301, 0.001574214227981301)
    13 = @(1.088159666666667, 3.264479, 0.002123718167154624, 0.002123718167154624, 3, 1, 10.38032666666667, 31.14098, 0.1781209566020688, 0.1781209566020688, 11.29545280549111, 101.65907524942, 0.000378278711486774, 0.000378278711486774)
    14 = @(1.088159666666667, 3.264479, 0.002123718167154624, 0.002123718167154624, 2, 0.6666666666666666, 0.422089, 1.266267, 0.007242825670663926, 0.007242825670663927, 0.4593002255436666, 4.133702029893001, 0.00001538172045832285, 0.00001538172045832286)
    15 = @(1.088159666666667, 3.264479, 0.002123718167154624, 0.002123718167154624, 3, 1, 4.133026333333333, 12.399079, 0.07092056230936286, 0.07092056230936288, 4.497392557204555, 40.476533014841, 0.0001506152866012154, 0.0001506152866012154)
    16 = @(1.088159666666667, 3.264479, 0.002123718167154624, 0.002123718167154624, 3, 1, 0.1434673333333333, 0.430402, 0.002461824128959449, 0.002461824128959449, 0.1561153656175555, 1.405038290558, 0.000005228220627010791, 0.000005228220627010791)
    17 = @(19.43698366666667, 58.310951, 0.03793439197579861, 0.03793439197579861, 3, 1, 43.19793166666667, 129.593795, 0.7412538312889448, 0.7412538312889448, 839.6374922387828, 7556.737430149045, 0.02811901338967733, 0.02811901338967733)
    18 = @(19.43698366666667, 58.310951, 0.03793439197579861, 0.03793439197579861, 3, 1, 10.38032666666667, 31.14098, 0.1781209566020688, 0.1781209566020688, 201.7622398746645, 1815.86015887198, 0.006756910186847093, 0.006756910186847093)
    19 = @(19.43698366666667, 58.310951, 0.03793439197579861, 0.03793439197579861, 2, 0.6666666666666666, 0.422089, 1.266267, 0.007242825670663926, 0.007242825670663927, 8.204136998879667, 73.83723298991701, 0.0002747521880033419, 0.0002747521880033419)
    20 = @(19.43698366666667, 58.310951, 0.03793439197579861, 0.03793439197579861, 3, 1, 4.133026333333333, 12.399079, 0.07092056230936286, 0.07092056230936288, 80.33356533490321, 723.002088014129, 0.00269032840978742, 0.002690328409787421)
    21 = @(19.43698366666667, 58.310951, 0.03793439197579861, 0.03793439197579861, 3, 1, 0.1434673333333333, 0.430402, 0.002461824128959449, 0.002461824128959449, 2.788572214700222, 25.097149932302, 0.00009338780148342675, 0.00009338780148342675)
    22 = @(0.451366, 1.354098, 0.0008809131633892397, 0.0008809131633892397, 3, 1, 43.19793166666667, 129.593795, 0.7412538312889448, 0.7412538312889448, 19.49807762465667, 175.48269862191, 0.0006529802573951382, 0.0006529802573951382)
    23 = @(0.451366, 1.354098, 0.0008809131633892397, 0.0008809131633892397, 3, 1, 10.38032666666667, 31.14098, 0.1781209566020688, 0.1781209566020688, 4.685326526226667, 42.16793873604, 0.000156909095346246, 0.000156909095346246)
    24 = @(0.451366, 1.354098, 0.0008809131633892397, 0.0008809131633892397, 2, 0.6666666666666666, 0.422089, 1.266267, 0.007242825670663926, 0.007242825670663927, 0.190516623574, 1.714649612166, 0.000006380300473421351, 0.000006380300473421352)
    25 = @(0.451366, 1.354098, 0.0008809131633892397, 0.0008809131633892397, 3, 1, 4.133026333333333, 12.399079, 0.07092056230936286, 0.07092056230936288, 1.865507563971333, 16.789568075742, 0.00006247485689328452, 0.00006247485689328453)
    26 = @(0.451366, 1.354098, 0.0008809131633892397, 0.0008809131633892397, 3, 1, 0.1434673333333333, 0.430402, 0.002461824128959449, 0.002461824128959449, 0.06475627637733333, 0.5828064873960001, 0.000002168653281149628, 0.000002168653281149628)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 7 + $i   # column G = 7
        $ws.Cells.Item([int]$row, $col).Value = $values[$i]
    }
}

Write-Output "done"